$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 148 - shifts rows 148:243 down to 149:244,
# carrying formatting (e.g. the date style on column D) along with them.
$ws.Rows.Item(148).Insert()

# Populate the newly inserted row 148 with the new record.
$ws.Range("A148").Value = 4
$ws.Range("B148").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C148").Value = "Los Lagos"
$ws.Range("D148").Value = 44529
$ws.Range("E148").Value = 10
$ws.Range("F148").Value = 100112023
$ws.Range("G148").Value = "Brócoli"
$ws.Range("H148").Value = "Sin especificar"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 500
$ws.Range("K148").Value = 1200
$ws.Range("L148").Value = 1200
$ws.Range("M148").Value = 1200
$ws.Range("N148").Value = "$/unidad"
$ws.Range("O148").Value = "Región del Maule"
$ws.Range("P148").Value = 1200
$ws.Range("Q148").Value = 1
$ws.Range("R148").Value = "Hortaliza"
